# This script applies the crypto price/volume update described by the diff.
# The D (Price) and E (Volume(1h)) columns hold TEXT values in the original
# workbook (they were inline strings, e.g. "510.03", "  +2.26%  ", "56.932.21").
# Excel auto-converts plain-looking decimal strings (e.g. "513.97") to numbers
# when assigned through .Value, so for those specific values we first force the
# cells NumberFormat to Text ("@") to preserve the original text representation
# (this also keeps trailing zeros such as "6.70" or "1.20" intact).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Row swap: row 17 (WrappedEther) <-> row 18 (Polkadot) become Polkadot / WrappedEther ---
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D17" "6.11"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.032.32"
$ws.Range("E18").Value = "  +1.85%  "

# --- Price / Volume updates ---
$ws.Range("D2").Value = "57.221.60"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "3.041.40"
$ws.Range("E3").Value = "  +2.46%  "
$ws.Range("E4").Value = "  -0.06%  "
Set-TextCell "D5" "513.97"
$ws.Range("E5").Value = "  +3.53%  "
Set-TextCell "D6" "141.36"
$ws.Range("E6").Value = "  +3.59%  "
Set-TextCell "D7" "0.999"
$ws.Range("E7").Value = "  -0.09%  "
Set-TextCell "D8" "0.438"
$ws.Range("E8").Value = "  +3.02%  "
Set-TextCell "D9" "7.19"
$ws.Range("E9").Value = "  -1.36%  "
Set-TextCell "D10" "0.109"
$ws.Range("E10").Value = "  +3.13%  "
Set-TextCell "D11" "0.377"
$ws.Range("E11").Value = "  +6.40%  "
$ws.Range("D12").Value = "3.556.95"
$ws.Range("E12").Value = "  +2.02%  "
Set-TextCell "D13" "0.126"
$ws.Range("E13").Value = "  -1.77%  "
Set-TextCell "D14" "26.91"
$ws.Range("E14").Value = "  +4.97%  "
Set-TextCell "D15" "0.0000166"
$ws.Range("E15").Value = "  +6.34%  "
$ws.Range("D16").Value = "57.197.56"
$ws.Range("E16").Value = "  +0.57%  "
Set-TextCell "D19" "13.35"
$ws.Range("E19").Value = "  +6.62%  "
Set-TextCell "D20" "8.08"
$ws.Range("E20").Value = "  +4.36%  "
Set-TextCell "D21" "333.75"
$ws.Range("E21").Value = "  +4.75%  "
$ws.Range("E22").Value = "  +0.19%  "
Set-TextCell "D23" "0.506"
$ws.Range("E23").Value = "  +4.75%  "
Set-TextCell "D24" "65.43"
$ws.Range("E24").Value = "  +3.57%  "
$ws.Range("D25").Value = "3.159.09"
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("E26").Value = "  +0.20%  "
Set-TextCell "D27" "0.164"
$ws.Range("E27").Value = "  +2.25%  "
$ws.Range("D28").Value = "0.0₃0922"
$ws.Range("E28").Value = "  +4.48%  "
Set-TextCell "D29" "6.70"
$ws.Range("E29").Value = "  +2.80%  "
Set-TextCell "D30" "7.16"
$ws.Range("E30").Value = "  +1.21%  "
Set-TextCell "D31" "1.81"
$ws.Range("E31").Value = "  +2.96%  "
Set-TextCell "D32" "1.20"
$ws.Range("E32").Value = "  +4.26%  "
Set-TextCell "D33" "20.69"
$ws.Range("E33").Value = "  +3.10%  "
Set-TextCell "D34" "4.69"
$ws.Range("E34").Value = "  +2.64%  "
Set-TextCell "D35" "153.80"
$ws.Range("E35").Value = "  +0.45%  "
Set-TextCell "D36" "5.94"
$ws.Range("E36").Value = "  +4.11%  "
Set-TextCell "D37" "1.28"
$ws.Range("E37").Value = "  +3.22%  "
Set-TextCell "D38" "25.44"
$ws.Range("E38").Value = "  +6.78%  "
Set-TextCell "D39" "0.0674"
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("D40").Value = "3.069.91"
$ws.Range("E40").Value = "  +2.40%  "
Set-TextCell "D41" "37.04"
$ws.Range("E41").Value = "  -1.01%  "
Set-TextCell "D42" "3.87"
$ws.Range("E42").Value = "  +4.82%  "
$ws.Range("E43").Value = "  -0.02%  "
Set-TextCell "D44" "0.662"
$ws.Range("E44").Value = "  +3.82%  "
$ws.Range("D45").Value = "2.207.10"
$ws.Range("E45").Value = "  +0.72%  "
Set-TextCell "D46" "1.39"
$ws.Range("E46").Value = "  +1.18%  "
Set-TextCell "D47" "0.958"
$ws.Range("E47").Value = "  +1.78%  "
Set-TextCell "D48" "6.02"
$ws.Range("E48").Value = "  +2.02%  "
Set-TextCell "D51" "0.0173"
$ws.Range("E51").Value = "  +14.29%  "

# --- Row swap: row 49 (InjectiveProtocol) <-> row 50 (VeChain) become VeChain / InjectiveProtocol ---
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D49" "0.0243"
$ws.Range("E49").Value = "  +4.20%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D50" "20.26"
$ws.Range("E50").Value = "  +6.94%  "
